$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Add new row of data (row 3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Customer MIA"
$ws.Range("C3").Value = "Resolved"
$ws.Range("D3").Value = "New customer acquired"
# Use the raw date serial number (2013-10-12) instead of a DateTime object so
# Excel doesn't auto-register a brand new number-format style for the cell;
# the date style is applied explicitly below via PasteSpecial from E2.
$ws.Range("E3").Value = 41559

# Copy the date formatting/style from E2 to E3
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)  # xlPasteFormats

# Widen column D (target stored width 21.1640625; engine quantizes ColumnWidth
# to the nearest 1/6 character unit on write, so 20.3 is the input that lands
# on the closest achievable stored width, 21.1666666...)
$ws.Columns.Item(4).ColumnWidth = 20.3

# Update selection to E4, matching diff
$ws.Range("E4").Select()
